$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move CPU_Table from K1:M3 to M1:O3 ---
$loCPU = $ws.ListObjects.Item("CPU_Table")
$loCPU.Resize($ws.Range("M1:O3"))

# --- Expand Process_Table from A1:I32 to A1:K32 ---
$loP = $ws.ListObjects.Item("Process_Table")
$loP.Resize($ws.Range("A1:K32"))

# Write new header/data for CPU_Table at its new location (M1:O3)
$ws.Cells.Item(1, 13).Value2 = "CPU ID:"
$ws.Cells.Item(1, 14).Value2 = "Busy Time (ms):"
$ws.Cells.Item(1, 15).Value2 = "Idle Time (ms):"
$ws.Cells.Item(3, 13).Value2 = "Average"
$ws.Cells.Item(3, 14).Formula = "=AVERAGE(CPU_Table[Busy Time (ms):])"
$ws.Cells.Item(3, 15).Formula = "=AVERAGE(CPU_Table[Idle Time (ms):])"

# Clear old CPU_Table cells that are not reused by Process_Table (L1, M1, L3, K3 - columns 12,13,11)
$ws.Cells.Item(1, 12).ClearContents()
$ws.Cells.Item(3, 12).ClearContents()
$ws.Cells.Item(3, 11).ClearContents()

# New headers for Process_Table
$ws.Cells.Item(1, 10).Value2 = "# of Page Faults:"
$ws.Cells.Item(1, 11).Value2 = "Average page fault servicing time (ns):"

# Wrap text style on K1 (new header cell)
$ws.Cells.Item(1, 11).WrapText = $true

# Totals row formulas
$ws.Cells.Item(32, 10).Formula = "=AVERAGE(Process_Table['# of Page Faults:])"
$ws.Cells.Item(32, 11).Formula = "=AVERAGE(Process_Table[Average page fault servicing time (ns):])"

# Row height for header row (to accommodate wrapped text)
$ws.Rows.Item(1).RowHeight = 30

# Column widths
$ws.Columns.Item(9).ColumnWidth = 12.85546875
$ws.Columns.Item(10).ColumnWidth = 17.7109375
$ws.Columns.Item(11).ColumnWidth = 20.5703125
$ws.Columns.Item(13).ColumnWidth = 9.85546875
$ws.Columns.Item(14).ColumnWidth = 17.42578125
$ws.Columns.Item(15).ColumnWidth = 16.7109375

# Selection
$ws.Range("K2").Select()
